# Update gh-pages to output generated at 456a3b4
# Increments the "F" column (view/heat count) for several rows that are
# duplicated across the "展览" sheet and the combined "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F9").Value  = 152
$ws1.Range("F10").Value = 152
$ws1.Range("F12").Value = 3065
$ws1.Range("F26").Value = 4726
$ws1.Range("F29").Value = 1624
$ws1.Range("F30").Value = 65
$ws1.Range("F31").Value = 114

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F17").Value = 152
$ws4.Range("F18").Value = 152
$ws4.Range("F21").Value = 3065
$ws4.Range("F37").Value = 4726
$ws4.Range("F40").Value = 1624
$ws4.Range("F43").Value = 65
$ws4.Range("F44").Value = 114
